# Add two new columns "I0" (col I) and "IF" (col J) to the stats sheet,
# mirroring the existing header style used by the other header cells
# (bold font, thin border, centered/top-aligned) and filling in the
# per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header row (B1:H1): bold, thin box
# border, centered horizontally, top-aligned vertically.
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("I1:J1").VerticalAlignment = -4160     # xlTop
$ws.Range("I1:J1").Borders.LineStyle = 1         # xlContinuous

# --- Data ------------------------------------------------------------------
$data = @{
    2  = @(8, 8)
    3  = @(5, 6)
    4  = @(8, 8)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(7, 8)
    8  = @(9, 9)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(7, 8)
    15 = @(8, 8)
    16 = @(6, 7)
    17 = @(9, 9)
    18 = @(8, 9)
    19 = @(8, 8)
    20 = @(7, 7)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(8, 9)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(8, 8)
    29 = @(1, 4)
    30 = @(1, 3)
    31 = @(4, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value  = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]   # column J
}
